$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "2025/12/02 16:12"
$ws.Range("B2").Value = "-"
$ws.Range("C2").Value = "-"
$ws.Range("D2").Value = "-"
$ws.Range("E2").Value = "-"
$ws.Range("F2").Value = "-"
$ws.Range("G2").Value = "-"
